$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "66.577.26"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "3.241.04"
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue "D5" "605.30"
$ws.Range("E5").Value = "  +0.16%  "
Set-TextValue "D6" "156.57"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.240.50"
$ws.Range("E8").Value = "  +1.64%  "
Set-TextValue "D9" "0.548"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("E11").Value = "  -1.25%  "
Set-TextValue "D12" "0.501"
$ws.Range("E12").Value = "  -1.88%  "
Set-TextValue "D13" "0.0000273"
$ws.Range("E13").Value = "  +2.90%  "
Set-TextValue "D14" "39.07"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "3.767.16"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "66.613.24"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "3.251.92"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("E19").Value = "  +1.63%  "
Set-TextValue "D20" "508.09"
$ws.Range("E20").Value = "  -1.05%  "
Set-TextValue "D21" "15.32"
$ws.Range("E21").Value = "  -0.87%  "
Set-TextValue "D22" "0.743"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E23").Value = "  -1.44%  "
Set-TextValue "D24" "14.64"
$ws.Range("E24").Value = "  -1.81%  "
Set-TextValue "D25" "86.11"
$ws.Range("E25").Value = "  +1.69%  "
Set-TextValue "D26" "0.171"
$ws.Range("E26").Value = "  +90.58%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +0.12%  "
Set-TextValue "D29" "9.07"
$ws.Range("E29").Value = "  -1.57%  "
Set-TextValue "D30" "2.34"
$ws.Range("E30").Value = "  -2.26%  "
Set-TextValue "D31" "2.91"
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("E32").Value = "  -2.20%  "
Set-TextValue "D33" "28.19"
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -4.57%  "
Set-TextValue "D36" "6.35"
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("D37").Value = "0.0₃0815"
$ws.Range("E37").Value = "  +20.61%  "
Set-TextValue "D38" "55.37"
$ws.Range("E38").Value = "  +1.16%  "
Set-TextValue "D39" "495.52"
$ws.Range("E39").Value = "  -3.52%  "
Set-TextValue "D40" "3.25"
$ws.Range("E40").Value = "  +13.68%  "
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("E42").Value = "  +2.16%  "
Set-TextValue "D43" "8.74"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("E44").Value = "  -3.70%  "
$ws.Range("D45").Value = "2.947.52"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("E46").Value = "  +0.91%  "
Set-TextValue "D47" "28.17"
$ws.Range("E47").Value = "  -0.83%  "
Set-TextValue "D48" "2.40"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("E50").Value = "  -0.04%  "
Set-TextValue "D51" "2.55"
$ws.Range("E51").Value = "  -1.80%  "
